# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp (A1)
# - Update case statistics for several countries (Pakistan, Armenia/Kirguistan,
#   Australia, Tailandia, Butan, Montserrat/Islas Malvinas)
# - Kirguistan now overtakes Armenia in the ranking (row 55), and
#   Islas Malvinas overtakes Montserrat (row 213)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 13 de Agosto de 2020 a las 07:18"

# Pakistan (row 17)
$ws.Range("B17").Value = 286674
$ws.Range("C17").Value = 753
$ws.Range("D17").Value = 264060
$ws.Range("E17").Value = 16475
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 6139

# Row 55: was Armenia, now Kirguistan (overtakes Armenia in ranking)
$ws.Range("A55").Value = "Kirguistan"
$ws.Range("B55").Value = 41069
$ws.Range("C55").Value = 310
$ws.Range("D55").Value = 33288
$ws.Range("E55").Value = 6294
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 3
$ws.Range("H55").Value = 1487

# Row 56: was Kirguistan, now Armenia
$ws.Range("A56").Value = "Armenia"
$ws.Range("B56").Value = 40794
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 33492
$ws.Range("E56").Value = 6496
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 806

# Australia (row 72)
$ws.Range("B72").Value = 22358
$ws.Range("C72").Value = 231
$ws.Range("D72").Value = 12779
$ws.Range("E72").Value = 9218
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 9
$ws.Range("H72").Value = 361

# Tailandia (row 117)
$ws.Range("B117").Value = 3359
$ws.Range("C117").Value = 3
$ws.Range("D117").Value = 3169
$ws.Range("E117").Value = 132
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 58

# Butan (row 192)
$ws.Range("B192").Value = 116
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 97
$ws.Range("E192").Value = 19
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 0

# Row 213: was Montserrat, now Islas Malvinas (overtakes Montserrat)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214: was Islas Malvinas, now Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
